# Cord Cutting Wizard Change Log - Aug 2019 update
#
# Re-word the "Comment" category labels used in column E. The old combined
# comment "Network added to Service (base or Add-On Package) in Aug 2019"
# is split into two distinct comments depending on whether the network was
# added to the base service (column D == "Yes") or to a named Add-On
# Package.
#
# The category labels are applied in a specific pass order so that the
# workbook's shared-string table ends up with the new label text in the
# same relative order as the published workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 336

function Update-Comments($oldText, $newText) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        if ($cell.Value() -eq $oldText) {
            $cell.Value = $newText
        }
    }
}

# 1) Alias Changed for Network
Update-Comments "New Alias added for Network name in Aug 2019" "Alias Changed for Network"

# 2) Name of Add-On Package Changed
Update-Comments "Name of Add-On Package changed in Aug 2019" "Name of Add-On Package Changed"

# 3) Network Added to Base Service (split of the old combined comment, base case)
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value() -eq "Network added to Service (base or Add-On Package) in Aug 2019") {
        $newValueCell = $ws.Cells.Item($r, 4)
        if ($newValueCell.Value() -eq "Yes") {
            $cell.Value = "Network Added to Base Service"
        }
    }
}

# 4) Network Removed from Database
Update-Comments "Old Network removed from database in Aug 2019" "Network Removed from Database"

# 5) Network Moved from Base Service to Add-On Package
Update-Comments "Network moved within Service from Base to Add-On Package" "Network Moved from Base Service to Add-On Package"

# 6) New Network Added to Database in Aug 2019
Update-Comments "New Network added to database in Aug 2019" "New Network Added to Database in Aug 2019"

# 7) Network Removed from Add-On Package
Update-Comments "Network removed from Service (base or Add-on Package) in Aug 2019" "Network Removed from Add-On Package"

# 8) Network Added to Add-On Package (split of the old combined comment, add-on case)
Update-Comments "Network added to Service (base or Add-On Package) in Aug 2019" "Network Added to Add-On Package"

# Restore the active cell/selection on the (only) worksheet to A2, matching
# the saved view state of the updated workbook.
$ws.Range("A2").Select()
